$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset has shrunk from 10 rows to 6 rows, so remove the now-unused
# trailing rows first.
$ws.Rows("7:10").Delete()

# Update the remaining 6 rows with the refreshed "cloud" data set.
$ws.Range("A1").Value = "Torrontés Clásico"
$ws.Range("B1").Value = 1150
$ws.Range("C1").Value = "Bodega Dos"
$ws.Range("D1").Value = "Reconocida por sus Malbecs"
$ws.Range("E1").Value = "Argentina"
$ws.Range("F1").Value = "Torrontés ClásicoMourvèdre es una variedad de uva tinta que se utiliza en la producción de vinos tintos robustos y especiados.Cabernet Sauvignon es una variedad de uva tinta ampliamente reconocida por su presencia en los vinos tintos de Bordeaux.Garnacha Blanca es una variedad de uva blanca que produce vinos blancos con cuerpo y textura, con sabores a frutas blancas y notas florales."
$ws.Range("G1").Value = 9
$ws.Range("A2").Value = "Merlot Reserva"
$ws.Range("B2").Value = 1250
$ws.Range("C2").Value = "Bodega Tres"
$ws.Range("D2").Value = "Región importante de San Juan"
$ws.Range("E2").Value = "Argentina"
$ws.Range("F2").Value = "Merlot ReservaMerlot Reserva"
$ws.Range("G2").Value = 8.9
$ws.Range("A3").Value = "Chardonnay"
$ws.Range("B3").Value = 1350
$ws.Range("C3").Value = "Bodega Cuatro"
$ws.Range("D3").Value = "Famoso por sus vinos tintos de alta calidad"
$ws.Range("E3").Value = "Argentina"
$ws.Range("F3").Value = "Sémillon es una variedad de uva blanca que se utiliza en la producción de vinos blancos secos, dulces y también vinos de postre.Pinot Grigio es una variedad de uva blanca que produce vinos blancos ligeros y refrescantes, con notas cítricas y florales.Garnacha Blanca es una variedad de uva blanca que produce vinos blancos con cuerpo y textura, con sabores a frutas blancas y notas florales."
$ws.Range("G3").Value = 8.4
$ws.Range("A4").Value = "Pinot Noir"
$ws.Range("B4").Value = 1500
$ws.Range("C4").Value = "Bodega Tres"
$ws.Range("D4").Value = "Región importante de San Juan"
$ws.Range("E4").Value = "Argentina"
$ws.Range("F4").Value = "Garnacha Blanca es una variedad de uva blanca que produce vinos blancos con cuerpo y textura, con sabores a frutas blancas y notas florales.Tannat es una variedad de uva tinta que se asocia principalmente con los vinos de Uruguay, conocidos por su estructura tánica y sabor intenso."
$ws.Range("G4").Value = 8.25
$ws.Range("A5").Value = "Cabernet Sauvignon"
$ws.Range("B5").Value = 1300
$ws.Range("C5").Value = "Bodega Uno"
$ws.Range("D5").Value = "Famoso por sus vinos tintos de alta calidad"
$ws.Range("E5").Value = "Argentina"
$ws.Range("F5").Value = "Bonarda es una variedad de uva tinta que se utiliza en la producción de vinos tintos suaves y afrutados, con sabores a frutas negras y especias.Cabernet Sauvignon es una variedad de uva tinta ampliamente reconocida por su presencia en los vinos tintos de Bordeaux."
$ws.Range("G5").Value = 7.9
$ws.Range("A6").Value = "Malbec Reserva"
$ws.Range("B6").Value = 1200
$ws.Range("C6").Value = "Bodega Uno"
$ws.Range("D6").Value = "Famoso por sus vinos tintos de alta calidad"
$ws.Range("E6").Value = "Argentina"
$ws.Range("F6").Value = "Mourvèdre es una variedad de uva tinta que se utiliza en la producción de vinos tintos robustos y especiados.Grenache es una variedad de uva tinta que se utiliza en muchos vinos tintos y rosados, conocidos por su cuerpo y sabor afrutado."
$ws.Range("G6").Value = 6
